$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "71.531.93"
Set-TextValue "E2" "  +3.21%  "

Set-TextValue "D3" "3.717.25"
Set-TextValue "E3" "  +8.26%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "585.51"
Set-TextValue "E5" "  +0.93%  "

Set-TextValue "D6" "180.04"

Set-TextValue "D7" "3.715.05"
Set-TextValue "E7" "  +8.42%  "

Set-TextValue "E8" "  +3.97%  "

Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  -0.03%  "

Set-TextValue "E10" "  +2.29%  "

Set-TextValue "D11" "0.613"
Set-TextValue "E11" "  +4.64%  "

Set-TextValue "D12" "49.59"
Set-TextValue "E12" "  +1.88%  "

Set-TextValue "E13" "  +2.71%  "

Set-TextValue "D14" "4.315.75"
Set-TextValue "E14" "  +8.80%  "

Set-TextValue "D15" "680.36"
Set-TextValue "E15" "  -3.07%  "

Set-TextValue "D16" "9.07"
Set-TextValue "E16" "  +4.95%  "

Set-TextValue "D17" "3.692.85"
Set-TextValue "E17" "  +7.80%  "

Set-TextValue "D18" "71.689.73"
Set-TextValue "E18" "  +3.26%  "

Set-TextValue "E19" "  +1.29%  "

Set-TextValue "D20" "18.08"
Set-TextValue "E20" "  +1.91%  "

Set-TextValue "D21" "11.66"
Set-TextValue "E21" "  +1.97%  "

Set-TextValue "D22" "6.47"
Set-TextValue "E22" "  +20.06%  "

Set-TextValue "E23" "  +5.09%  "

Set-TextValue "D24" "17.51"
Set-TextValue "E24" "  +3.28%  "

Set-TextValue "D25" "102.79"
Set-TextValue "E25" "  +1.51%  "

Set-TextValue "E26" "  +3.12%  "

Set-TextValue "E27" "  +6.45%  "

Set-TextValue "D28" "10.48"
Set-TextValue "E28" "  +9.05%  "

Set-TextValue "D29" "35.72"
Set-TextValue "E29" "  +6.18%  "

Set-TextValue "D30" "9.23"
Set-TextValue "E30" "  +5.38%  "

Set-TextValue "E31" "  +5.79%  "

Set-TextValue "E32" "  +10.39%  "

Set-TextValue "D33" "587.71"
Set-TextValue "E33" "  +3.27%  "

Set-TextValue "D34" "11.26"
Set-TextValue "E34" "  +2.04%  "

Set-TextValue "E35" "  +3.70%  "

Set-TextValue "D36" "59.23"
Set-TextValue "E36" "  +1.95%  "

Set-TextValue "E37" "  +0.09%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D38" "3.692.60"
Set-TextValue "E38" "  +2.41%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.146"
Set-TextValue "E39" "  +5.46%  "

Set-TextValue "D40" "0.0₃0782"
Set-TextValue "E40" "  +6.75%  "

Set-TextValue "D41" "35.75"
Set-TextValue "E41" "  +2.40%  "

Set-TextValue "D42" "3.46"
Set-TextValue "E42" "  +5.53%  "

Set-TextValue "E43" "  +4.79%  "

Set-TextValue "D44" "0.0461"
Set-TextValue "E44" "  +9.65%  "

Set-TextValue "E45" "  +5.11%  "

Set-TextValue "E46" "  +8.67%  "

Set-TextValue "D47" "3.38"
Set-TextValue "E47" "  +0.88%  "

Set-TextValue "E48" "  +4.26%  "

Set-TextValue "D49" "1.46"
Set-TextValue "E49" "  -0.72%  "

Set-TextValue "D50" "0.998"

Set-TextValue "D51" "136.29"
Set-TextValue "E51" "  +3.69%  "
